$d = $word.ActiveDocument

$pairs = @(
    @("897÷8=112, 1", "427÷7=61, 0"),
    @("610÷3=203, 1", "968÷4=242, 0"),
    @("337÷8=42, 1", "809÷9=89, 8"),
    @("703÷3=234, 1", "968÷4=242, 0"),
    @("822÷5=164, 2", "850÷5=170, 0"),
    @("512÷5=102, 2", "891÷2=445, 1"),
    @("651÷7=93, 0", "851÷4=212, 3"),
    @("195÷5=39, 0", "662÷2=331, 0"),
    @("711÷8=88, 7", "195÷9=21, 6"),
    @("994÷9=110, 4", "328÷3=109, 1"),
    @("825÷2=412, 1", "143÷6=23, 5"),
    @("848÷2=424, 0", "151÷4=37, 3"),
    @("163÷5=32, 3", "188÷7=26, 6"),
    @("121÷6=20, 1", "850÷3=283, 1"),
    @("558÷3=186, 0", "317÷4=79, 1"),
    @("287÷4=71, 3", "950÷9=105, 5"),
    @("226÷8=28, 2", "727÷6=121, 1"),
    @("377÷6=62, 5", "820÷3=273, 1"),
    @("519÷8=64, 7", "725÷8=90, 5"),
    @("391÷7=55, 6", "876÷3=292, 0"),
    @("637÷2=318, 1", "896÷9=99, 5"),
    @("844÷3=281, 1", "760÷8=95, 0"),
    @("680÷2=340, 0", "230÷6=38, 2"),
    @("238÷8=29, 6", "337÷9=37, 4"),
    @("366÷8=45, 6", "302÷7=43, 1")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
